# Update the "Newest" form number counters on row 3 of Sheet 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1030
$ws.Range("C3").Value = 2003
$ws.Range("D3").Value = 3003
$ws.Range("E3").Value = 4003
